# Regenerate the participant's task-order workbook: each worksheet keeps its
# underlying physical identity (and therefore its r:id in workbook.xml.rels),
# but tab names are re-assigned and the generated stimulus-file listings are
# replaced with a freshly "randomized" batch (new run timestamps), matching
# the new "experiment order generation script".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (physical sheet1.xml): was GNG (4 rows) -> becomes TOL (6 rows)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TOL_TO-16515889613168297"

$ws1.Range("B2").Value = "MM_stims-16515889612848053.csv"
$ws1.Range("B3").Value = "ZM_stims-16515889612671185.csv"
$ws1.Range("B4").Value = "MM_stims-16515889613003638.csv"
$ws1.Range("B5").Value = "ZM_stims-16515889612858076.csv"

# Two new rows - copy the existing styled cell down first so the new "A"
# cells pick up the same bold/border/center style (s="1") as the rest of
# the column, then overwrite the values.
$ws1.Range("A5").Copy($ws1.Range("A6"))
$ws1.Range("A5").Copy($ws1.Range("A7"))
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "MM_stims-16515889613155239.csv"
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "ZM_stims-16515889613014126.csv"

# ---------------------------------------------------------------------
# Sheet 2 (physical sheet2.xml): was NB (9 rows) -> stays NB (9 rows)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651588963511486"

$ws2.Range("B2").Value = "ZB-match_8-16515889617887733.csv"
$ws2.Range("B3").Value = "TB-16515889625981922.csv"
$ws2.Range("B4").Value = "OB-16515889620207105.csv"
$ws2.Range("B5").Value = "ZB-match_5-16515889617588763.csv"
$ws2.Range("B6").Value = "TB-16515889624283602.csv"
$ws2.Range("B7").Value = "TB-16515889634937243.csv"
$ws2.Range("B8").Value = "OB-16515889618930802.csv"
$ws2.Range("B9").Value = "ZB-match_0-16515889615588672.csv"
$ws2.Range("B10").Value = "OB-16515889618211083.csv"

# ---------------------------------------------------------------------
# Sheet 3 (physical sheet3.xml): was RS (2 rows) -> becomes vSAT (4 rows)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "vSAT_TO-16515889635910282"

$ws3.Range("B2").Value = "vSAT_stims-16515889635749226.csv"
$ws3.Range("B3").Value = "vSAT_stims-16515889635590837.csv"

$ws3.Range("A3").Copy($ws3.Range("A4"))
$ws3.Range("A3").Copy($ws3.Range("A5"))
$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "SAT_stims-1651588963521913.csv"
$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "SAT_stims-16515889635425766.csv"

# ---------------------------------------------------------------------
# Sheet 4 (physical sheet4.xml): was TOL (6 rows) -> becomes GNG (4 rows)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "GNG_TO-16515889636246362"

$ws4.Range("B2").Value = "go_stims-16515889635957766.csv"
$ws4.Range("B3").Value = "GNG_stims-1651588963606868.csv"
$ws4.Range("B4").Value = "go_stims-16515889636098762.csv"
$ws4.Range("B5").Value = "GNG_stims-16515889636226845.csv"

# Drop the two rows TOL had that GNG does not - fully clear (value + style)
# so the cells disappear from the sheet and the used range shrinks back to
# A1:B5, matching the smaller GNG table.
$ws4.Range("A6:B7").Clear()

# ---------------------------------------------------------------------
# Sheet 5 (physical sheet5.xml): was vSAT (4 rows) -> becomes RS (2 rows)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "RS_TO-16515889636257298"

$ws5.Range("B2").Value = "eyes closed"
$ws5.Range("B3").Value = "eyes open"

# Drop the two rows vSAT had that RS does not.
$ws5.Range("A4:B5").Clear()

# NOTE: tab order itself does not change here - both before and after, the
# Nth tab is backed by the Nth physical worksheet part (rId/sheetN.xml), so
# no .Move() calls are needed; only names + cell contents were rewritten.
